$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Could not find text: $old"
    }
}

# 1. "...based on recent results in the photometry and thermodynamics of radiation literature."
#    -> "...based on recent results in the fields of radiation thermodynamics and photometry."
Replace-Text "based on recent results in the photometry and thermodynamics of radiation literature." "based on recent results in the fields of radiation thermodynamics and photometry."

# 2 & 3. Rewrite of "The exact method (a) is free of ... (b) ... and (c) allows choices for the spectral weighting function."
Replace-Text "is free of (possibly erroneous) assumptions for the maximum luminous efficacy, (b) uses the correct spectral exergy-to-energy ratio, no longer assuming its value to be 1, and" "is free of any assumptions for the value of the maximum luminous efficacy, (b)"

Replace-Text "(c) allows choices for the spectral weighting function." "uses a non-unity spectral exergy-to-energy ratio, and (c) allows choices for the spectral luminous weighting function, which converts broad-spectrum electromagnetic radiation to light."

# 4. "...assumptions inherent to the original method and leads..." -> "...assumptions inherent to the conventional method and leads..."
Replace-Text "assumptions inherent to the original method and leads" "assumptions inherent to the conventional method and leads"

# 5. "...growing body of societal exergy analysis literature, and it concludes with specific recommendations for societal exergy analysts. "
#    -> "...growing field of societal exergy analysis, and it concludes with specific recommendations for societal exergy practitioners. "
Replace-Text "is important, because it clarifies a methodological issue in the growing body of societal exergy analysis literature, and it concludes with specific recommendations for societal exergy analysts. " "is important, because it clarifies a methodological issue in the growing field of societal exergy analysis, and it concludes with specific recommendations for societal exergy practitioners. "

# 6. ". In particular, the following are directly relevant to our article" -> ". In particular, the following subject areas are directly relevant to our article"
Replace-Text ". In particular, the following are directly relevant to our article" ". In particular, the following subject areas are directly relevant to our article"

# 7. "find the manuscript to be both novel and important " -> "agree that the manuscript is both novel and important "
Replace-Text "find the manuscript to be both novel and important " "agree that the manuscript is both novel and important "
